# Updated symbol list on Tue Dec 27 09:36:36 UTC 2022 with GitHub Actions
# Refreshes the "Price" (column D) and "Volume(1h)" label (column E) cells
# with the latest scraped values for the crypto ranking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prices (column D) are stored as text, so values that look numeric are
# written with a leading apostrophe to force a text interpretation and
# keep exact formatting (e.g. trailing zeros like "0.0006000").
$ws.Range("D2").Value  = "'242.70"
$ws.Range("D3").Value  = "'22.99"
$ws.Range("D4").Value  = "'5.395"
$ws.Range("D5").Value  = "'0.05946"
$ws.Range("D6").Value  = "'3.428"
$ws.Range("D7").Value  = "'6.501"
$ws.Range("D8").Value  = "'0.8137"
$ws.Range("D9").Value  = "'0.9285"
$ws.Range("D10").Value = "'0.1430"
$ws.Range("D11").Value = "'0.07425"
$ws.Range("D12").Value = "'0.03281"
$ws.Range("D13").Value = "'0.03085"
$ws.Range("D14").Value = "'0.09351"
$ws.Range("D15").Value = "'3.865"
$ws.Range("D16").Value = "'0.001574"
$ws.Range("D17").Value = "'0.04699"

$ws.Range("D18").Value = "'0.0006000"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("D19").Value = "'0.005891"

$ws.Range("D20").Value = "'0.001258"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("D21").Value = "'0.004788"
$ws.Range("D22").Value = "'0.00008000"
$ws.Range("D24").Value = "'2.160"
$ws.Range("D26").Value = "'0.1330"

$ws.Range("E27").Value = "26UpBotsUBXTWorstin24h"

$ws.Range("D41").Value = "'0.006350"

$ws.Range("D42").Value = "'0.003800"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"

$ws.Range("D44").Value = "'0.008905"
$ws.Range("D45").Value = "'0.00005160"
$ws.Range("D47").Value = "'0.7000"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.0002000"
